$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") date values for rows 2-5 from 45224 (2023-10-25)
# to 45233 (2023-11-03), matching the diff.
$ws.Range("C2").Value = 45233
$ws.Range("C3").Value = 45233
$ws.Range("C4").Value = 45233
$ws.Range("C5").Value = 45233
